# Revert the "acknowledgement of service" -> "acknowledgement of claim" rename.
# The heading currently reads "Acknowledgement" + " of " + "claim" (three runs).
# Replacing the whole-word "claim" run with "service" merges its text into the
# preceding " of " run (producing " of service") and removes the now-empty
# "claim" run, leaving the "Acknowledgement" run untouched - exactly matching
# the target markup. MatchCase + MatchWholeWord keep this from touching other
# occurrences of "claim" used as template/merge-field text (e.g. "claimIssueDate").
$d = $word.ActiveDocument
$d.Content.Find.Execute("claim", $true, $true, $false, $false, $false,
                         $true, 1, $false, "service", 2)
